$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value2 = "(TerrierRetr(BM25) >> (TerrierRetr(Tf) ** TerrierRetr(PL2) ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score()) >> <pyterrier.ltr.LTRTransformer object at 0x000001A33AE4CAD0>)"
$ws.Range("B4").Value2 = 0.2636832138379338
$ws.Range("C4").Value2 = 0.4328666394185419

$ws.Range("A5").Value2 = "(<pyterrier._ops.SetUnion object at 0x000001949B4976F0> >> <pyterrier.terrier._text_loader.TerrierTextLoader object at 0x00000194A6A5FC50> >> (<pyterrier.terrier.retriever.TextScorer object at 0x000001949B4975C0> ** <pyterrier.terrier.retriever.TextScorer object at 0x000001949B496EA0> ** <pyterrier.terrier.retriever.TextScorer object at 0x0000019496D74A70> ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score()) >> <pyterrier.ltr.LTRTransformer object at 0x00000194A6A5F110>)"
$ws.Range("B5").Value2 = 0.09141840205407331
$ws.Range("C5").Value2 = 0.3209493560511296

$ws.Range("A6").Value2 = "(<pyterrier._ops.SetUnion object at 0x0000019251400050> >> <pyterrier.terrier._text_loader.TerrierTextLoader object at 0x00000192514B97F0> >> (<pyterrier.terrier.retriever.TextScorer object at 0x00000192514B9550> ** <pyterrier.terrier.retriever.TextScorer object at 0x0000019251400190> ** <pyterrier.terrier.retriever.TextScorer object at 0x00000192514002D0> ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score() ** pt.apply.doc_score()) >> <pyterrier.ltr.LTRTransformer object at 0x00000192514B9A90>)"
$ws.Range("B6").Value2 = 0.2743828363576117
$ws.Range("C6").Value2 = 0.4760306064989383
